$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Notes": update the issue-description note (A3)
# ---------------------------------------------------------------------------
$wsNotes = $wb.Worksheets.Item("Notes")
$wsNotes.Range("A3").Value = "Specific issue: survey_ids (in surveys table) are not unique"

# ---------------------------------------------------------------------------
# Sheet "studies": full header/data restructure
# ---------------------------------------------------------------------------
$wsStudies = $wb.Worksheets.Item("studies")
$wsStudies.Cells.Clear()

$wsStudies.Range("A1").Value = "study_id"
$wsStudies.Range("B1").Value = "study_label"
$wsStudies.Range("C1").Value = "description"
$wsStudies.Range("D1").Value = "access_level"
$wsStudies.Range("E1").Value = "contributors"
$wsStudies.Range("F1").Value = "reference"
$wsStudies.Range("G1").Value = "reference_year"

$wsStudies.Range("A2").Value = "foo"
$wsStudies.Range("D2").Value = "public"
$wsStudies.Range("F2").Value = "https://doi.org/10.1093%2Fgenetics%2F16.2.97"
$wsStudies.Range("F2").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Sheet "surveys": full header/data restructure (13 columns now)
# ---------------------------------------------------------------------------
$wsSurveys = $wb.Worksheets.Item("surveys")
$wsSurveys.Cells.Clear()

$wsSurveys.Range("A1").Value = "study_id"
$wsSurveys.Range("B1").Value = "survey_id"
$wsSurveys.Range("C1").Value = "country_name"
$wsSurveys.Range("D1").Value = "site_name"
$wsSurveys.Range("E1").Value = "latitude"
$wsSurveys.Range("F1").Value = "longitude"
$wsSurveys.Range("G1").Value = "location_method"
$wsSurveys.Range("H1").Value = "location_notes"

$wsSurveys.Range("I1").NumberFormat = "@"
$wsSurveys.Range("I1").Value = "collection_start"
$wsSurveys.Range("J1").NumberFormat = "@"
$wsSurveys.Range("J1").Value = "collection_end"
$wsSurveys.Range("K1").NumberFormat = "@"
$wsSurveys.Range("K1").Value = "collection_day"
$wsSurveys.Range("L1").NumberFormat = "@"
$wsSurveys.Range("L1").Value = "time_method"

$wsSurveys.Range("M1").Value = "time_notes"

# row 2
$wsSurveys.Range("A2").Value = "foo"
$wsSurveys.Range("B2").Value = "S01"
$wsSurveys.Range("E2").Value = 0
$wsSurveys.Range("F2").Value = 0
$wsSurveys.Range("H2").Value = "example data"
$wsSurveys.Range("K2").NumberFormat = "@"
$wsSurveys.Range("K2").Value = "2020-01-01"
$wsSurveys.Range("L2").NumberFormat = "@"
$wsSurveys.Range("M2").Value = "example data"

# row 3
$wsSurveys.Range("A3").Value = "foo"
$wsSurveys.Range("B3").Value = "S01"
$wsSurveys.Range("E3").Value = 0
$wsSurveys.Range("F3").Value = 0
$wsSurveys.Range("H3").Value = "example data"
$wsSurveys.Range("K3").NumberFormat = "@"
$wsSurveys.Range("K3").Value = "2020-01-01"
$wsSurveys.Range("L3").NumberFormat = "@"
$wsSurveys.Range("M3").Value = "example data"

# ---------------------------------------------------------------------------
# Sheet "counts": rename two headers, update one data value
# ---------------------------------------------------------------------------
$wsCounts = $wb.Worksheets.Item("counts")
$wsCounts.Range("A1").Value = "study_id"
$wsCounts.Range("B1").Value = "survey_id"
$wsCounts.Range("A2").Value = "foo"

# ---------------------------------------------------------------------------
# Selections / active sheet & cell (matches new sheetViews in the diff)
# ---------------------------------------------------------------------------
$wsNotes.Activate()
$wsNotes.Range("A4").Select()

$wsStudies.Activate()
$wsStudies.Range("D3").Select()

$wsSurveys.Activate()
$wsSurveys.Range("C8").Select()

$wsCounts.Activate()
$wsCounts.Range("D7").Select()

# "studies" tab is the active one when the file is reopened
$wsStudies.Activate()
